# Fill in "reserva_total.prn" (Sheet1) with the "Analisis de la Reserva Total" report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- small helpers ------------------------------------------------------
# Style "header" rows: thin border all around + centered horizontal alignment
function Set-HeaderBoxStyle($rng) {
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
}
# Style "sub-header" rows: centered horizontal alignment only, no border
function Set-CenterStyle($rng) {
    $rng.HorizontalAlignment = -4108
}

# ---- Row 1: title ---------------------------------------------------------
$ws.Range("A1").Value = "Análisis de la Reserva Total"

# ---- Row 3: section header (border + center) ------------------------------
$ws.Range("A3:F3").Merge()
Set-HeaderBoxStyle $ws.Range("A3:F3")
$ws.Range("A3").Value = "RESERVA ROTANTE EN MAQUINAS QUE REGULAN"

# ---- Row 4 ----
$ws.Range("A4:C4").Merge()
$ws.Range("D4:F4").Merge()
$ws.Range("A4").Value = "RESERVA HIDRO [MW]"
$ws.Range("D4").Value = 195

# ---- Row 5 ----
$ws.Range("A5:C5").Merge()
$ws.Range("D5:F5").Merge()
$ws.Range("A5").Value = "RESERVA TERMICA [MW]"
$ws.Range("D5").Value = 515.4

# ---- Row 6 ----
$ws.Range("A6:C6").Merge()
$ws.Range("D6:F6").Merge()
$ws.Range("A6").Value = "RESERVA TOTAL [MW]"
$ws.Range("D6").Value = 710.4

# ---- Row 7: sub-header (center only, no border) - style only the anchor cell ----
$ws.Range("A7:E7").Merge()
Set-CenterStyle $ws.Range("A7")
$ws.Range("A7").Value = "RESERVA ROTANTE DEL PARQUE REGULANTE [%]"
$ws.Range("F7").Value = 21.81

# ---- Row 8: section header (border + center) ----
$ws.Range("A8:F8").Merge()
Set-HeaderBoxStyle $ws.Range("A8:F8")
$ws.Range("A8").Value = "RESERVA PROGRAMADA A 50Hz PARA RPF"

# ---- Row 9 ----
$ws.Range("A9:C9").Merge()
$ws.Range("D9:F9").Merge()
$ws.Range("A9").Value = "RESERVA HIDRO [MW]"
$ws.Range("D9").Value = 37.5

# ---- Row 10 ----
$ws.Range("A10:C10").Merge()
$ws.Range("D10:F10").Merge()
$ws.Range("A10").Value = "RESERVA TÉRMICA [MW]"
$ws.Range("D10").Value = 135.5

# ---- Row 11 ----
$ws.Range("A11:C11").Merge()
$ws.Range("D11:F11").Merge()
$ws.Range("A11").Value = "TOTAL SISTEMA [MW]"
$ws.Range("D11").Value = 173

# ---- Row 12: merged label, no special style ----
$ws.Range("A12:E12").Merge()
$ws.Range("A12").Value = "RESERVA PARA RPF [%]"
$ws.Range("F12").Value = 5.31

# ---- Row 13 ----
$ws.Range("A13:E13").Merge()
$ws.Range("A13").Value = "COLABORACIÓN DEL PARQUE HIDRO EN RSF [MW]"
$ws.Range("F13").Value = 157.5

# ---- Row 14 ----
$ws.Range("A14:E14").Merge()
$ws.Range("A14").Value = "COLABORACIÓN DEL PARQUE HIDRO EN RSF [%]"
$ws.Range("F14").Value = 4.84

# ---- Row 15: sub-header (center only, no border) - style only the anchor cell ----
$ws.Range("A15:F15").Merge()
Set-CenterStyle $ws.Range("A15")
$ws.Range("A15").Value = "POTENCIA OPERABLE EN EL PARQUE REGULANTE"

# ---- Row 16 ----
$ws.Range("A16:C16").Merge()
$ws.Range("D16:F16").Merge()
$ws.Range("A16").Value = "HIDRO [MW]"
$ws.Range("D16").Value = 945

# ---- Row 17 ----
$ws.Range("A17:C17").Merge()
$ws.Range("D17:F17").Merge()
$ws.Range("A17").Value = "TÉRMICA TG-CC [MW]"
$ws.Range("D17").Value = 1745.4

# ---- Row 18 ----
$ws.Range("A18:C18").Merge()
$ws.Range("D18:F18").Merge()
$ws.Range("A18").Value = "TÉRMICA TV [MW]"
$ws.Range("D18").Value = 900

# ---- Row 19 ----
$ws.Range("A19:C19").Merge()
$ws.Range("D19:F19").Merge()
$ws.Range("A19").Value = "TOTAL [MW]"
$ws.Range("D19").Value = 3590.4

# ---- Row 20: section header (border + center) ----
$ws.Range("A20:F20").Merge()
Set-HeaderBoxStyle $ws.Range("A20:F20")
$ws.Range("A20").Value = "RESERVA PROGRAMADA EN EL PARQUE REGULANTE"

# ---- Row 21 ----
$ws.Range("A21:C21").Merge()
$ws.Range("D21:F21").Merge()
$ws.Range("A21").Value = "HIDRO"
$ws.Range("D21").Value = 100

# ---- Row 22 ----
$ws.Range("A22:C22").Merge()
$ws.Range("D22:F22").Merge()
$ws.Range("A22").Value = "TÉRMICA TG-CC"
$ws.Range("D22").Value = 415.4

# ---- Row 23 ----
$ws.Range("A23:C23").Merge()
$ws.Range("D23:F23").Merge()
$ws.Range("A23").Value = "TÉRMICA TV"
$ws.Range("D23").Value = 100

# ---- Row 24 ----
$ws.Range("A24:C24").Merge()
$ws.Range("D24:F24").Merge()
$ws.Range("A24").Value = "TOTAL"
$ws.Range("D24").Value = 515.4

# ---- Row 25: not merged ----
$ws.Range("A25").Value = "RESERVA NUEVA"
$ws.Range("D25").Value = 488.5525634765625

# ---- Row 26: not merged ----
$ws.Range("A26").Value = "RESERVA TOTAL 2"
$ws.Range("D26").Value = 610.4
